$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 55: update note text, hours, and row height ---
$ws.Range("E55").Value = "Compiled all old data with new data (excluding overlap). Running and testing old code with new dataset, everything working fine"
$ws.Range("C55").Value = 6
$ws.Rows.Item(55).RowHeight = 45

# --- Row 56: add Hours + Notes, growing the row height ---
$ws.Range("E55").Copy()
$ws.Range("E56").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C56").Value = 8
$ws.Range("E56").Value = "implementing residency, gender, citizenship. Have to rewrite and recompile new rows into old df to include gender etc."
$ws.Rows.Item(56).RowHeight = 45

# --- Row 57: brand new row of data ---
$ws.Range("A56").Copy()
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("B56").Copy()
$ws.Range("B57").PasteSpecial(-4122)
$ws.Range("E55").Copy()
$ws.Range("E57").PasteSpecial(-4122)

$ws.Range("A57").Value = 45495
$ws.Range("B57").Value = "M"
$ws.Range("C57").Value = 4
$ws.Range("E57").Value = "added new rows to both studentgrades and studentgrades_prof. new columns for new variables. Doesn't seem to be helping, in fact performing worse… need to debug"
$ws.Rows.Item(57).RowHeight = 60

$excel.CutCopyMode = 0

# --- Update the view so the new row / cell is the active selection ---
$ws.Range("C57").Select()
$excel.ActiveWindow.ScrollRow = 49
